$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (border/alignment/bold) for new rows 9:15 in column A
# by copying the existing formatted cell A2 down to A9:A15 (format only).
$ws.Range("A2:A8").Copy()
$ws.Range("A9:A15").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# Update the index column A (row number - 2) and ticker columns B-F per row.
# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "NSE:5PAISA"
$ws.Range("D2").Value = "NSE:EICHERMOT"
$ws.Range("F2").Value = "NSE:DIXON"

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "NSE:ATUL"
$ws.Range("F3").Value = "NSE:EICHERMOT"

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "NSE:AVANTIFEED"

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "NSE:BECTORFOOD"

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "NSE:DIXON"

# Row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "NSE:DODLA"

# Row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "NSE:GSPL"

# Row 9
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "NSE:JUSTDIAL"

# Row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "NSE:KRIDHANINF"

# Row 11
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "NSE:METROPOLIS"

# Row 12
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "NSE:MONTECARLO"

# Row 13
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "NSE:RALLIS"

# Row 14
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NSE:RATNAMANI"

# Row 15
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "NSE:RPPL"

# Clear out the old "support Zone" ticker values that moved out of column C
$ws.Range("C2").Value = ""
$ws.Range("C3").Value = ""
$ws.Range("C4").Value = ""
$ws.Range("C5").Value = ""
$ws.Range("C6").Value = ""
$ws.Range("C7").Value = ""
$ws.Range("C8").Value = ""
